$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "44-29=",
    "47+22=",
    "93-78=",
    "66-56=",
    "54-18=",
    "80-55=",
    "42+39=",
    "73-40=",
    "57-24=",
    "84-54=",
    "45+49=",
    "49+7=",
    "82-53=",
    "52+11=",
    "16+39=",
    "81-72=",
    "3+63=",
    "45+16=",
    "30+35=",
    "65-42=",
    "17+18=",
    "15-11=",
    "34-7=",
    "88-60=",
    "12+46=",
    "99-89=",
    "40+2=",
    "34+31=",
    "94-81=",
    "58-48=",
    "24+24=",
    "9+58=",
    "56-3=",
    "6+56=",
    "88-67=",
    "81-19=",
    "91+6=",
    "93-70=",
    "84-51=",
    "86-5=",
    "49+21=",
    "99-83=",
    "96-17=",
    "58+6=",
    "54-4=",
    "96-29=",
    "82-67=",
    "48-32=",
    "86-83=",
    "81-3=",
    "86-68=",
    "56-12=",
    "63-27=",
    "61+38=",
    "5+19=",
    "61-55=",
    "54-31=",
    "15+14=",
    "42+3=",
    "3+23=",
    "84-36=",
    "48-6=",
    "32+12=",
    "98-39=",
    "42+44=",
    "28+40=",
    "0+98=",
    "73-52=",
    "90-11=",
    "49+13=",
    "93-56=",
    "25+8=",
    "51+4=",
    "18-1=",
    "24+57=",
    "66-61=",
    "8-2=",
    "56+30=",
    "92-81=",
    "47-3=",
    "80-45=",
    "6+43=",
    "54+7=",
    "89-87=",
    "39-38=",
    "58-36=",
    "29+3=",
    "54+22=",
    "72-59=",
    "79-46=",
    "40+26=",
    "64-41=",
    "84-51=",
    "32+62=",
    "31+56=",
    "34+37=",
    "44+26=",
    "25-15=",
    "3-0=",
    "1+87="
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Output "Updated $idx cells"